$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.436.93"
$ws.Range("E2").Value = "  +2.22%  "

$ws.Range("D3").Value = "1.796.43"
$ws.Range("E3").Value = "  +3.12%  "

$ws.Range("E4").Value = "  +0.64%  "

$ws.Range("D5").Value = "337.24"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").Value = "0.3806"
$ws.Range("E7").Value = "  +1.51%  "

$ws.Range("D8").Value = "0.3454"
$ws.Range("E8").Value = "  +1.43%  "

$ws.Range("D9").Value = "48.38"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").Value = "1.202"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("D11").Value = "0.07488"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("D13").Value = "22.02"
$ws.Range("E13").Value = "  +7.04%  "

$ws.Range("D14").Value = "6.468"
$ws.Range("E14").Value = "  +0.51%  "

$ws.Range("D15").Value = "1.795.59"
$ws.Range("E15").Value = "  +3.52%  "

$ws.Range("D16").Value = "7.060"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  +1.55%  "

$ws.Range("D18").Value = "0.06643"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").Value = "84.80"
$ws.Range("E19").Value = "  +2.35%  "

$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").Value = "6.496"
$ws.Range("E21").Value = "  +4.36%  "

$ws.Range("D22").Value = "17.32"
$ws.Range("E22").Value = "  +3.33%  "

$ws.Range("D23").Value = "27.417.74"
$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("D24").Value = "12.50"
$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("D25").Value = "2.426"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "1.497"
$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.558"
$ws.Range("E27").Value = "  +4.76%  "

$ws.Range("D28").Value = "21.39"
$ws.Range("E28").Value = "  +9.06%  "

$ws.Range("D29").Value = "152.09"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "1.999.51"
$ws.Range("E30").Value = "  +3.59%  "

$ws.Range("D31").Value = "133.81"
$ws.Range("E31").Value = "  +0.98%  "

$ws.Range("D32").Value = "4.063"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").Value = "6.109"
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").Value = "0.08693"
$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("D35").Value = "13.24"
$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("D36").Value = "1.686"
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("D37").Value = "5.450"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "0.6888"
$ws.Range("E38").Value = "  +9.47%  "

$ws.Range("D39").Value = "8.908"
$ws.Range("E39").Value = "  +4.71%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06368"
$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.2206"
$ws.Range("E41").Value = "  +1.02%  "

$ws.Range("D42").Value = "0.02335"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("E43").Value = "  +4.20%  "

$ws.Range("D44").Value = "14.42"
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("D45").Value = "0.6437"
$ws.Range("E45").Value = "  +5.22%  "

$ws.Range("D46").Value = "0.9990"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("D47").Value = "3.863"
$ws.Range("E47").Value = "  -1.72%  "

$ws.Range("D48").Value = "2.125"
$ws.Range("E48").Value = "  +2.12%  "

$ws.Range("D49").Value = "129.87"
$ws.Range("E49").Value = "  +0.28%  "

$ws.Range("E50").Value = "  -0.37%  "

$ws.Range("D51").Value = "79.34"
$ws.Range("E51").Value = "  +1.65%  "
